# Update cryptos list data (prices & 1h volume %) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.708.50"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").Value = "'1.804.52"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'231.71"
$ws.Range("E5").Value = "  -1.86%  "

$ws.Range("D6").Value = "'0.5941"
$ws.Range("E6").Value = "  -2.57%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'0.2788"
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").Value = "'0.06836"

$ws.Range("D10").Value = "'23.35"
$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("D11").Value = "'0.07544"
$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.809.18"
$ws.Range("E12").Value = "  -0.88%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.806"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("D14").Value = "'0.6233"
$ws.Range("E14").Value = "  -1.20%  "

$ws.Range("D15").Value = "'2.049.86"
$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("D16").Value = "'0.000009332"
$ws.Range("E16").Value = "  -6.61%  "

$ws.Range("D17").Value = "'75.72"
$ws.Range("E17").Value = "  -3.58%  "

$ws.Range("D18").Value = "'28.655.22"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("E19").Value = "  -6.15%  "

$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "'210.56"
$ws.Range("E21").Value = "  -7.06%  "

$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("D23").Value = "'6.870"
$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").Value = "'154.31"
$ws.Range("E25").Value = "  -0.84%  "

$ws.Range("D26").Value = "'7.877"
$ws.Range("E26").Value = "  -2.00%  "

$ws.Range("E27").Value = "  -2.47%  "

$ws.Range("D28").Value = "'16.40"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").Value = "'1.426"
$ws.Range("E29").Value = "  -4.59%  "

$ws.Range("D30").Value = "'0.06183"
$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").Value = "'3.787"
$ws.Range("E32").Value = "  -0.88%  "

$ws.Range("D33").Value = "'3.756"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("E34").Value = "  -0.94%  "

$ws.Range("D35").Value = "'1.067"
$ws.Range("E35").Value = "  -4.94%  "

$ws.Range("D36").Value = "'0.6418"
$ws.Range("E36").Value = "  -0.35%  "

$ws.Range("E37").Value = "  -2.14%  "

$ws.Range("D38").Value = "'2.715"
$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("D39").Value = "'6.488"
$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("D40").Value = "'0.01715"
$ws.Range("E40").Value = "  -1.30%  "

$ws.Range("D41").Value = "'1.134.77"
$ws.Range("E41").Value = "  -6.62%  "

$ws.Range("D42").Value = "'0.8799"
$ws.Range("E42").Value = "  -3.54%  "

$ws.Range("D43").Value = "'1.007"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").Value = "'100.72"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "'1.967.18"
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").Value = "'60.56"
$ws.Range("E46").Value = "  -3.26%  "

$ws.Range("E47").Value = "  -2.70%  "

$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").Value = "'8.378"
$ws.Range("E49").Value = "  -2.08%  "

$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").Value = "'0.4486"
$ws.Range("E51").Value = "  -1.81%  "
